# Updated symbol list on Sun Jan  8 13:43:10 UTC 2023 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) snapshot values for the coin rows.
# Values are written as literal text (matching the source sheet, where every
# cell is stored as a plain string) rather than letting Excel reinterpret
# them as numbers/percentages, so number-like strings such as "0.9091" or
# "3.449" keep their exact original formatting/trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "262.03";    "E2"  = "0.43%"
    "D3"  = "26.63";     "E3"  = "-2.09%"
    "D4"  = "4.701";     "E4"  = "0.20%"
    "E5"  = "-0.66%"
    "D6"  = "6.701";     "E6"  = "0.68%"
    "E7"  = "-0.28%"
    "D8"  = "0.9091";    "E8"  = "-1.25%"
    "D9"  = "0.1402";    "E9"  = "-0.18%"
    "D10" = "0.05099";   "E10" = "8.95%"
    "E11" = "0.28%"
    "D12" = "0.03112";   "E12" = "1.76%"
    "E13" = "-0.25%"
    "D14" = "0.001529";  "E14" = "0.25%"
    "D15" = "0.0006188"; "E15" = "1.72%"
    "D16" = "0.005969";  "E16" = "-1.19%"
    "D17" = "3.449";     "E17" = "-0.05%"
    "D18" = "3.166";     "E18" = "0.65%"
    "D19" = "2.147";     "E19" = "-0.78%"
    "D21" = "0.1281";    "E21" = "-2.23%"
    "D22" = "4.129";     "E22" = "1.12%"
    "D23" = "0.04235";   "E23" = "-0.15%"
    "D25" = "0.004062";  "E25" = "6.91%"
    "E26" = "0.06%"
    "E27" = "23.06%"
    "D40" = "0.03955";   "E40" = "2.20%"
    "D41" = "0.1111";    "E41" = "-0.07%"
    "D42" = "0.004192";  "E42" = "2.84%"
    "D43" = "0.01393";   "E43" = "-14.75%"
    "E44" = "-7.01%"
    "D45" = "0.00005126"; "E45" = "-0.56%"
    "E46" = "0.06%"
    "D48" = "0.2579";    "E48" = "58.99%"
    "E49" = "0.06%"
    "E50" = "0.06%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
